# "added new room and tests"
# Extends the maze grid on Sheet1 with a new room "J" to the right of the
# existing wall column T (rows 9-16), including its doorway cells "JL".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Values -----------------------------------------------------------
$ws.Range("T9:T16").Value = "w"

$ws.Range("U9").Value = "w"
$ws.Range("V9").Value = "w"
$ws.Range("W9").Value = "w"

$ws.Range("U10:W11").Value = "J"
$ws.Range("V12:W13").Value = "J"
$ws.Range("U14:W15").Value = "J"

$ws.Range("U12").Value = "JL"
$ws.Range("U13").Value = "JL"

$ws.Range("U16").Value = "w"
$ws.Range("V16").Value = "w"

# ---- Formatting (reuse existing fills instead of inventing new ones) --
$xlPasteFormats = -4122

# "w" wall formatting - copy from an existing wall cell's format.
$ws.Range("W9").Copy()
$ws.Range("T9:T16").PasteSpecial($xlPasteFormats)
$ws.Range("U9:W9").PasteSpecial($xlPasteFormats)
$ws.Range("U16:V16").PasteSpecial($xlPasteFormats)

# "J" room formatting - copy from an existing room ("M") cell's format.
$ws.Range("Q10").Copy()
$ws.Range("U10:W11").PasteSpecial($xlPasteFormats)
$ws.Range("V12:W13").PasteSpecial($xlPasteFormats)
$ws.Range("U14:W15").PasteSpecial($xlPasteFormats)

# "JL" doorway formatting - copy from an existing doorway cell's format.
$ws.Range("S9").Copy()
$ws.Range("U12:U13").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ---- View / selection state saved with the workbook -------------------
$ws.Range("U10").Select()
